$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.662.60'
$ws.Range('E2').Value = '  -1.87%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.797.89'
$ws.Range('E3').Value = '  -1.57%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '231.41'
$ws.Range('E5').Value = '  -1.35%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5890'
$ws.Range('E6').Value = '  -1.92%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2768'
$ws.Range('E8').Value = '  -0.36%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06800'
$ws.Range('E9').Value = '  -3.18%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '23.22'
$ws.Range('E10').Value = '  -0.58%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07532'
$ws.Range('E11').Value = '  -1.52%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.796.81'
$ws.Range('E12').Value = '  -1.71%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.782'
$ws.Range('E13').Value = '  -0.04%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6193'
$ws.Range('E14').Value = '  -0.61%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.042.52'
$ws.Range('E15').Value = '  -1.56%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.000009109'
$ws.Range('E16').Value = '  -7.87%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '75.56'
$ws.Range('E17').Value = '  -4.08%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '28.645.19'
$ws.Range('E18').Value = '  -2.00%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.478'
$ws.Range('E19').Value = '  -5.82%  '
$ws.Range('E20').Value = '  +0.06%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '210.96'
$ws.Range('E21').Value = '  -5.28%  '
$ws.Range('E22').Value = '  -1.09%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.822'
$ws.Range('E23').Value = '  -2.23%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.004'
$ws.Range('E24').Value = '  -0.05%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '153.68'
$ws.Range('E25').Value = '  -1.18%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.943'
$ws.Range('E26').Value = '  -0.15%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1269'
$ws.Range('E27').Value = '  -1.49%  '
$ws.Range('E28').Value = '  -0.15%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.423'
$ws.Range('E29').Value = '  -3.69%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.06131'
$ws.Range('E30').Value = '  -1.25%  '
$ws.Range('E31').Value = '  -1.24%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.824'
$ws.Range('E32').Value = '  +1.04%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.783'
$ws.Range('E33').Value = '  -0.84%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.734'
$ws.Range('E34').Value = '  -0.30%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.054'
$ws.Range('E35').Value = '  -5.05%  '
$ws.Range('E36').Value = '  +0.17%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.499'
$ws.Range('E37').Value = '  -1.64%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.715'
$ws.Range('E38').Value = '  -0.66%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.529'
$ws.Range('E39').Value = '  -0.12%  '
$ws.Range('E40').Value = '  -1.72%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.146.05'
$ws.Range('E41').Value = '  -5.98%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.8874'
$ws.Range('E42').Value = '  -0.94%  '
$ws.Range('E43').Value = '  +0.30%  '
$ws.Range('E44').Value = '  +0.02%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.950.43'
$ws.Range('E45').Value = '  -1.71%  '
$ws.Range('E46').Value = '  -3.21%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.00000000111'
$ws.Range('E47').Value = '  -3.46%  '
$ws.Range('E48').Value = '  +1.55%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.357'
$ws.Range('E49').Value = '  -1.80%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05471'
$ws.Range('E50').Value = '  -0.43%  '
$ws.Range('E51').Value = '  -1.65%  '
